$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value (values are text in the original sheet, so a leading
# apostrophe is used to force text interpretation and avoid numeric/percent coercion).
$updates = @{
    "D2" = "'304.67"
    "E2" = "'0.86%"
    "E3" = "'-3.82%"
    "D4" = "'5.098"
    "E4" = "'1.84%"
    "E5" = "'-0.17%"
    "D6" = "'2.152"
    "E6" = "'-2.90%"
    "D7" = "'7.934"
    "E7" = "'-1.13%"
    "B8" = "'GateToken"
    "C8" = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D8" = "'4.105"
    "E8" = "'2.10%"
    "B9" = "'MXToken"
    "C9" = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D9" = "'0.9192"
    "E9" = "'1.51%"
    "B10" = "'LiechtensteinCryptoassetsExchange"
    "C10" = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D10" = "'0.09712"
    "E10" = "'0.63%"
    "B11" = "'WazirX"
    "C11" = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D11" = "'0.1864"
    "E11" = "'-1.30%"
    "B12" = "'MandalaExchangeToken"
    "C12" = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D12" = "'0.08581"
    "E12" = "'1.25%"
    "B13" = "'BitrueCoin"
    "C13" = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D13" = "'0.03502"
    "E13" = "'-0.72%"
    "B14" = "'BitMartToken"
    "C14" = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D14" = "'0.09932"
    "E14" = "'-0.27%"
    "B15" = "'BitForexToken"
    "C15" = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D15" = "'0.001451"
    "E15" = "'-2.63%"
    "B16" = "'TigerCash"
    "C16" = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D16" = "'0.005696"
    "E16" = "'0.63%"
    "B17" = "'LEO"
    "C17" = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D17" = "'3.463"
    "E17" = "'-0.03%"
    "D18" = "'2.488"
    "E18" = "'20.46%"
    "E19" = "'-1.04%"
    "D20" = "'0.1331"
    "E20" = "'1.64%"
    "D21" = "'4.767"
    "E21" = "'0.09%"
    "D22" = "'0.2205"
    "E22" = "'0.00%"
    "D23" = "'0.04547"
    "E23" = "'-2.18%"
    "D24" = "'0.005082"
    "E24" = "'14.18%"
    "D25" = "'0.001237"
    "E25" = "'0.69%"
    "D27" = "'0.0004760"
    "E27" = "'0.19%"
    "D39" = "'0.01840"
    "E39" = "'4.70%"
    "D40" = "'0.04739"
    "E40" = "'0.54%"
    "D41" = "'0.007652"
    "E41" = "'-1.48%"
    "D42" = "'0.1399"
    "E42" = "'0.61%"
    "D43" = "'0.007743"
    "E43" = "'0.99%"
    "E44" = "'3.32%"
    "D45" = "'0.01118"
    "E45" = "'13.28%"
    "D46" = "'0.00006399"
    "E46" = "'5.07%"
    "E47" = "'0.20%"
    "D48" = "'0.0005813"
    "E48" = "'0.21%"
    "D49" = "'42.11"
    "E49" = "'385.65%"
    "E50" = "'-25.50%"
    "D51" = "'0.00002105"
    "E51" = "'0.20%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
